$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8549507260322571
$ws.Range("B1").Value = 2.643798351287842
$ws.Range("C1").Value = 3.292206048965454
$ws.Range("D1").Value = 1.852099180221558
$ws.Range("E1").Value = 1.418121337890625
